# Add a new "severity" column (INFO/WARNING/ERROR/CRITICAL) to both the
# Korean ("한국어") and English sheets, inserted right before the existing
# "actionLabel" column (old column F). This shifts the old F/G/H
# (actionLabel/actionType/actionTarget) columns one place to the right,
# becoming G/H/I.

$wb = $excel.ActiveWorkbook

# severity values per data row (row 2..6), same on both sheets
$severities = @{
    2 = "WARNING"
    3 = "ERROR"
    4 = "INFO"
    5 = "ERROR"
    6 = "CRITICAL"
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert a new blank column at F; existing F/G/H (and their widths)
    # shift right to G/H/I automatically, and the dimension/refs update.
    $ws.Columns("F:F").Insert()

    # Give the freshly inserted column a width in line with its neighbours
    # (the sheet's other "default" columns are ~12.8 characters wide).
    $ws.Columns("F:F").ColumnWidth = 12

    # Header
    $ws.Range("F1").Value = "severity"

    # Data rows
    foreach ($row in $severities.Keys) {
        $ws.Range("F$row").Value = $severities[$row]
    }
}
